# Update PPT with new images
#
# - Remove the two trailing duplicate slides (old slide 3 and slide 4).
# - On the remaining first slide, fix the first picture's alt text and
#   reposition/resize both pictures so they sit with a margin instead of
#   filling the whole slide.
# - On the remaining second slide (old slide 2), swap the alt text of its
#   two pictures (so it now matches what used to be slide 3's captions)
#   and apply the same repositioning/resizing.

$p = $ppt.ActivePresentation

# EMU-per-point conversion used throughout (1 pt = 12700 EMU).
$emuPerPt = 12700

# --- Remove the two extra slides -----------------------------------------
# Delete from the end first so earlier indices stay valid.
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()

# --- Slide 1 (was slide 1): fix alt text + reposition pictures -----------
$s1 = $p.Slides.Item(1)

$pic1 = $s1.Shapes.Item(1)
$pic1.AlternativeText = "image.jpg"
$pic1.Left = 457200 / $emuPerPt
$pic1.Top = 457200 / $emuPerPt
$pic1.Width = 3886200 / $emuPerPt
$pic1.Height = 5943600 / $emuPerPt

$pic2 = $s1.Shapes.Item(2)
$pic2.Left = 4800600 / $emuPerPt
$pic2.Top = 457200 / $emuPerPt
$pic2.Width = 3886200 / $emuPerPt
$pic2.Height = 5943600 / $emuPerPt

# --- Slide 2 (was slide 2): swap alt text + reposition pictures ----------
$s2 = $p.Slides.Item(2)

$s2pic1 = $s2.Shapes.Item(1)
$s2pic2 = $s2.Shapes.Item(2)

$tmpAlt = $s2pic1.AlternativeText
$s2pic1.AlternativeText = $s2pic2.AlternativeText
$s2pic2.AlternativeText = $tmpAlt

$s2pic1.Left = 457200 / $emuPerPt
$s2pic1.Top = 457200 / $emuPerPt
$s2pic1.Width = 3886200 / $emuPerPt
$s2pic1.Height = 5943600 / $emuPerPt

$s2pic2.Left = 4800600 / $emuPerPt
$s2pic2.Top = 457200 / $emuPerPt
$s2pic2.Width = 3886200 / $emuPerPt
$s2pic2.Height = 5943600 / $emuPerPt
